# ---------------------------------------------------------------------------
# Applies three edits to the manuscript:
#   1. Bumps the title-page date from 2022-05-10 to 2022-05-11.
#   2. Inserts "therefore" into the AR1 process sentence in the methods text.
#   3. Replaces the "Describe rho and phi" TODO bullet with the full
#      explanatory paragraph (with inline rho/phi/sigma equations) that
#      describes how temporal/species correlation and process-error
#      magnitude are controlled.
# ---------------------------------------------------------------------------

$d = $word.ActiveDocument

# --- Change 1: update the date on the title page -------------------------
$d.Content.Find.Execute("2022-05-10", $true, $false, $false, $false, $false, $true, 1, $false, "2022-05-11", 2) | Out-Null

# --- Change 2: add "therefore" to the AR1 sentence ------------------------
$oldSentence = "A first-order autoregressive (AR1) process was applied to account for temporal dependence."
$newSentence = "A first-order autoregressive (AR1) process was therefore applied to account for temporal dependence."
$d.Content.Find.Execute($oldSentence, $true, $false, $false, $false, $false, $true, 1, $false, $newSentence, 2) | Out-Null

# --- Change 3: replace the "Describe rho and phi" placeholder bullet with
#     the full paragraph of explanatory text (incl. OMath equations). ------
$target = $null
$paras = $d.Paragraphs
for ($i = 1; $i -le $paras.Count; $i++) {
    $p = $paras.Item($i)
    if ($p.Range.Text -like "Describe rho and phi*") {
        $target = $p
        break
    }
}

if ($target -ne $null) {
    $xml = '<?xml version="1.0" encoding="UTF-8" standalone="yes"?><pkg:package xmlns:pkg="http://schemas.microsoft.com/office/2006/xmlPackage"><pkg:part pkg:name="/word/document.xml" pkg:contentType="application/vnd.openxmlformats-officedocument.wordprocessingml.document.main+xml"><pkg:xmlData><w:document xmlns:w="http://schemas.openxmlformats.org/wordprocessingml/2006/main" xmlns:m="http://schemas.openxmlformats.org/officeDocument/2006/math"><w:body><w:p><w:pPr><w:pStyle w:val="FirstParagraph"/></w:pPr><w:r><w:t xml:space="preserve">The degree of temporal correlation is controlled by</w:t></w:r><w:r><w:t xml:space="preserve"> </w:t></w:r><m:oMath><m:r><m:t>ϕ</m:t></m:r></m:oMath><w:r><w:t xml:space="preserve">, where low to high correlation is represented by values between 0 and 1, and species-to-species correlations are described by</w:t></w:r><w:r><w:t xml:space="preserve"> </w:t></w:r><m:oMath><m:sSub><m:e><m:r><m:t>ρ</m:t></m:r></m:e><m:sub><m:r><m:t>s</m:t></m:r><m:r><m:rPr><m:sty m:val="p"/></m:rPr><m:t>,</m:t></m:r><m:r><m:t>s</m:t></m:r></m:sub></m:sSub></m:oMath><w:r><w:t xml:space="preserve">, where negative an positive correlation is represented by values between -1 and 1. This is a flexible structure that allows for the testing of alternate hypotheses that process errors are independent through time or across species (i.e,</w:t></w:r><w:r><w:t xml:space="preserve"> </w:t></w:r><m:oMath><m:r><m:t>ϕ</m:t></m:r><m:r><m:rPr><m:sty m:val="p"/></m:rPr><m:t>=</m:t></m:r><m:r><m:t>0</m:t></m:r></m:oMath><w:r><w:t xml:space="preserve"> </w:t></w:r><w:r><w:t xml:space="preserve">or</w:t></w:r><w:r><w:t xml:space="preserve"> </w:t></w:r><m:oMath><m:sSub><m:e><m:r><m:t>ρ</m:t></m:r></m:e><m:sub><m:r><m:t>s</m:t></m:r><m:r><m:rPr><m:sty m:val="p"/></m:rPr><m:t>,</m:t></m:r><m:r><m:t>s</m:t></m:r></m:sub></m:sSub><m:r><m:rPr><m:sty m:val="p"/></m:rPr><m:t>=</m:t></m:r><m:r><m:t>0</m:t></m:r></m:oMath><w:r><w:t xml:space="preserve">). The possibility that process errors are similarly correlated across all species may also be tested by estimating only one</w:t></w:r><w:r><w:t xml:space="preserve"> </w:t></w:r><m:oMath><m:r><m:t>ρ</m:t></m:r></m:oMath><w:r><w:t xml:space="preserve"> </w:t></w:r><w:r><w:t xml:space="preserve">parameter. Finally, the magnitude of the process error deviations are controlled by the species-specific standard deviation parameters,</w:t></w:r><w:r><w:t xml:space="preserve"> </w:t></w:r><m:oMath><m:sSub><m:e><m:r><m:t>σ</m:t></m:r></m:e><m:sub><m:r><m:t>s</m:t></m:r></m:sub></m:sSub></m:oMath><w:r><w:t xml:space="preserve">.</w:t></w:r></w:p></w:body></w:document></pkg:xmlData></pkg:part></pkg:package>'
    $target.Range.InsertXML($xml) | Out-Null
}
